# Weekly update: insert a new "Ajo" (garlic) price record for the
# "Mercado Mayorista Lo Valledor de Santiago" market as row 257, pushing
# the existing rows 257-347 down to 258-348 (dimension grows to A1:R348).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 257; everything below shifts down one row.
$ws.Rows("257:257").Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(257, 1).Value  = 6
$ws.Cells.Item(257, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(257, 3).Value  = "Metropolitana"
$ws.Cells.Item(257, 4).Value  = 44466
$ws.Cells.Item(257, 5).Value  = 13
$ws.Cells.Item(257, 6).Value  = 100112003
$ws.Cells.Item(257, 7).Value  = "Ajo"
$ws.Cells.Item(257, 8).Value  = "Chino"
$ws.Cells.Item(257, 9).Value  = "Primera"
$ws.Cells.Item(257, 10).Value = 1500
$ws.Cells.Item(257, 11).Value = 14500
$ws.Cells.Item(257, 12).Value = 15000
$ws.Cells.Item(257, 13).Value = 14700
$ws.Cells.Item(257, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(257, 15).Value = "China"
$ws.Cells.Item(257, 16).Value = 1470
$ws.Cells.Item(257, 17).Value = 10
$ws.Cells.Item(257, 18).Value = "Hortaliza"
